$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates as described by the diff (crypto price/volume refresh,
# plus two row swaps: Stellar/InternetComputer(DFINITY) and PaxDollar/TheSandbox).
# Column D values that look like plain numbers are written with a leading
# apostrophe (quote-prefix) so Excel stores them as text, matching the
# original inline-string cell type instead of silently coercing to a number.

$ws.Cells.Item(2, 4).Value = '30.277.88'
$ws.Cells.Item(2, 5).Value = '  -0.67%  '
$ws.Cells.Item(3, 4).Value = '1.858.31'
$ws.Cells.Item(3, 5).Value = '  -2.23%  '
$ws.Cells.Item(4, 4).Value = '''0.9999'
$ws.Cells.Item(5, 4).Value = '''236.07'
$ws.Cells.Item(5, 5).Value = '  -1.19%  '
$ws.Cells.Item(6, 4).Value = '''0.9998'
$ws.Cells.Item(7, 5).Value = '  -2.26%  '
$ws.Cells.Item(8, 4).Value = '''0.2805'
$ws.Cells.Item(8, 5).Value = '  -3.95%  '
$ws.Cells.Item(9, 4).Value = '''0.06477'
$ws.Cells.Item(9, 5).Value = '  -2.98%  '
$ws.Cells.Item(10, 4).Value = '1.854.62'
$ws.Cells.Item(10, 5).Value = '  -2.49%  '
$ws.Cells.Item(11, 4).Value = '''0.07382'
$ws.Cells.Item(12, 4).Value = '''16.28'
$ws.Cells.Item(12, 5).Value = '  -4.19%  '
$ws.Cells.Item(13, 4).Value = '''5.018'
$ws.Cells.Item(13, 5).Value = '  -3.40%  '
$ws.Cells.Item(14, 4).Value = '''87.17'
$ws.Cells.Item(14, 5).Value = '  -1.25%  '
$ws.Cells.Item(15, 4).Value = '''0.6451'
$ws.Cells.Item(15, 5).Value = '  -3.44%  '
$ws.Cells.Item(16, 4).Value = '30.221.92'
$ws.Cells.Item(16, 5).Value = '  -0.77%  '
$ws.Cells.Item(17, 5).Value = '  -0.02%  '
$ws.Cells.Item(18, 4).Value = '''13.14'
$ws.Cells.Item(18, 5).Value = '  -2.23%  '
$ws.Cells.Item(19, 4).Value = '''0.000007568'
$ws.Cells.Item(19, 5).Value = '  -3.76%  '
$ws.Cells.Item(20, 4).Value = '''223.22'
$ws.Cells.Item(20, 5).Value = '  +13.43%  '
$ws.Cells.Item(21, 4).Value = '2.095.23'
$ws.Cells.Item(21, 5).Value = '  -2.15%  '
$ws.Cells.Item(22, 4).Value = '''1.000'
$ws.Cells.Item(23, 4).Value = '''5.264'
$ws.Cells.Item(23, 5).Value = '  -3.52%  '
$ws.Cells.Item(24, 4).Value = '''6.066'
$ws.Cells.Item(24, 5).Value = '  -1.43%  '
$ws.Cells.Item(25, 4).Value = '''9.187'
$ws.Cells.Item(25, 5).Value = '  -2.98%  '
$ws.Cells.Item(26, 4).Value = '''163.38'
$ws.Cells.Item(26, 5).Value = '  +0.04%  '
$ws.Cells.Item(27, 4).Value = '''18.55'
$ws.Cells.Item(27, 5).Value = '  +0.81%  '
$ws.Cells.Item(28, 4).Value = '''1.924'
$ws.Cells.Item(28, 5).Value = '  -0.77%  '
$ws.Cells.Item(29, 4).Value = '''1.434'
$ws.Cells.Item(29, 5).Value = '  -2.59%  '
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(30, 4).Value = '''4.250'
$ws.Cells.Item(30, 5).Value = '  -1.88%  '
$ws.Cells.Item(31, 2).Value = 'Stellar'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(31, 4).Value = '''0.09202'
$ws.Cells.Item(31, 5).Value = '  +0.34%  '
$ws.Cells.Item(32, 4).Value = '''3.960'
$ws.Cells.Item(32, 5).Value = '  -4.22%  '
$ws.Cells.Item(33, 4).Value = '''0.04955'
$ws.Cells.Item(33, 5).Value = '  -4.22%  '
$ws.Cells.Item(34, 4).Value = '''1.143'
$ws.Cells.Item(34, 5).Value = '  +2.82%  '
$ws.Cells.Item(35, 4).Value = '''0.7250'
$ws.Cells.Item(35, 5).Value = '  -2.00%  '
$ws.Cells.Item(36, 4).Value = '''2.687'
$ws.Cells.Item(36, 5).Value = '  -1.53%  '
$ws.Cells.Item(37, 4).Value = '''0.01831'
$ws.Cells.Item(37, 5).Value = '  -0.77%  '
$ws.Cells.Item(38, 4).Value = '''2.596'
$ws.Cells.Item(38, 5).Value = '  -2.97%  '
$ws.Cells.Item(39, 4).Value = '''0.8969'
$ws.Cells.Item(39, 5).Value = '  -2.88%  '
$ws.Cells.Item(40, 4).Value = '''2.040'
$ws.Cells.Item(40, 5).Value = '  -1.26%  '
$ws.Cells.Item(41, 5).Value = '  -0.16%  '
$ws.Cells.Item(42, 4).Value = '''106.05'
$ws.Cells.Item(42, 5).Value = '  -1.04%  '
$ws.Cells.Item(43, 2).Value = 'TheSandbox'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(43, 4).Value = '''0.4245'
$ws.Cells.Item(43, 5).Value = '  -3.75%  '
$ws.Cells.Item(44, 2).Value = 'PaxDollar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(44, 4).Value = '''1.000'
$ws.Cells.Item(44, 5).Value = '  +0.48%  '
$ws.Cells.Item(45, 4).Value = '''7.270'
$ws.Cells.Item(45, 5).Value = '  -3.89%  '
$ws.Cells.Item(46, 5).Value = '  -5.78%  '
$ws.Cells.Item(47, 4).Value = '''63.31'
$ws.Cells.Item(47, 5).Value = '  -8.35%  '
$ws.Cells.Item(48, 4).Value = '''1.483'
$ws.Cells.Item(48, 5).Value = '  +5.83%  '
$ws.Cells.Item(49, 4).Value = '''8.707'
$ws.Cells.Item(49, 5).Value = '  -3.64%  '
$ws.Cells.Item(50, 4).Value = '''33.73'
$ws.Cells.Item(50, 5).Value = '  -3.53%  '
$ws.Cells.Item(51, 4).Value = '''0.05636'
$ws.Cells.Item(51, 5).Value = '  -3.37%  '
